$d = $word.ActiveDocument

# --- Skewness section: merge the proofing-split "<name>_p:-<number>" runs ---
# These three (degree_p, etest_p, mba_p) share the same "<word>_ + p : - + number"
# run layout, split apart by spellStart/gramStart/spellEnd/gramEnd proofErr marks.
# Re-finding "<name>_p" and ":-<number>" as contiguous text merges the runs Word
# split them into, while leaving the spellStart/spellEnd pair around the word intact.
$d.Content.Find.Execute("degree_p", $true, $false, $false, $false, $false, $true, 1, $false, "degree_p", 2) | Out-Null
$d.Content.Find.Execute(":-0.09749", $true, $false, $false, $false, $false, $true, 1, $false, ":-0.09749", 2) | Out-Null

$d.Content.Find.Execute("etest_p", $true, $false, $false, $false, $false, $true, 1, $false, "etest_p", 2) | Out-Null
$d.Content.Find.Execute(":-1.08858", $true, $false, $false, $false, $false, $true, 1, $false, ":-1.08858", 2) | Out-Null

$d.Content.Find.Execute("mba_p", $true, $false, $false, $false, $false, $true, 1, $false, "mba_p", 2) | Out-Null
$d.Content.Find.Execute(":-0.470723", $true, $false, $false, $false, $false, $true, 1, $false, ":-0.470723", 2) | Out-Null

# salary:-0.239837 is wrapped entirely in gramStart/gramEnd - re-finding the whole
# literal text merges it back into a single run.
$d.Content.Find.Execute("salary:-0.239837", $true, $false, $false, $false, $false, $true, 1, $false, "salary:-0.239837", 2) | Out-Null

# --- Kurtosis section: the three "Mesokurtic" explanations become "Platykurtic" ---
# Text changes from "Contains the positive value(s) falls under = 3" to
# "Contains the negative values falls under < 3" (does NOT touch the Leptokurtic
# paragraphs, which read "... falls under > 3 ...").
$d.Content.Find.Execute("Contains the positive value falls under = 3", $true, $false, $false, $false, $false, $true, 1, $false, "Contains the negative values falls under < 3", 2) | Out-Null

# The single combined replace above collapses "falls under < 3" into the
# surrounding non-bold run; restore bold on that phrase (also re-applies to the
# pre-existing Platykurtic paragraph harmlessly, since it is already bold there).
$rng = $d.Content
$rng.Find.ClearFormatting()
$found = $rng.Find.Execute("falls under < 3")
while ($found) {
  $rng.Bold = 1
  $rng.Collapse(0)
  $found = $rng.Find.Execute("falls under < 3")
}

# "Meso" + "kurtic" -> "Platykurtic" (keeps the surrounding bold parenthetical).
$d.Content.Find.Execute("Mesokurtic", $true, $false, $false, $false, $false, $true, 1, $false, "Platykurtic", 2) | Out-Null

Write-Output "done"
